$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tarantula")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.162617096018736
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 24.9268149882904

$ws = $wb.Worksheets.Item("Ochiai")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.216700819672132
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 23.02400468384074

$ws = $wb.Worksheets.Item("Op2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 3.128659250585479
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 24.44379391100702

$ws = $wb.Worksheets.Item("Barinel")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.162617096018736
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 24.9268149882904

$ws = $wb.Worksheets.Item("Dstar")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.200234192037472
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 23.07523419203747

$ws = $wb.Worksheets.Item("Russell_rao")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 5.122950819672131
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 27.92008196721311

$ws = $wb.Worksheets.Item("Simple_matching")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 26.88085480093671
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 75.98433840749414

$ws = $wb.Worksheets.Item("Rogers_tanimoto")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 26.88085480093671
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 75.98433840749414

$ws = $wb.Worksheets.Item("Ample")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.004464285714287
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 28.69767271662765

$ws = $wb.Worksheets.Item("Jaccard")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.90463992974239
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 22.74041276346604

$ws = $wb.Worksheets.Item("Cohen")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.90281030444965
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 42.69796545667447

$ws = $wb.Worksheets.Item("Scott")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 12.12492681498831
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 66.28732435597188

$ws = $wb.Worksheets.Item("Rogot1")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 12.12492681498831
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 66.28732435597188

$ws = $wb.Worksheets.Item("Geometric_mean")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.295374707259954
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 41.66239754098361

$ws = $wb.Worksheets.Item("M2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.063012295081968
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 24.20228337236534

$ws = $wb.Worksheets.Item("Wong1")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 5.122950819672131
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 27.92008196721311

$ws = $wb.Worksheets.Item("Sokal")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 26.88085480093671
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 75.98433840749414

$ws = $wb.Worksheets.Item("Sorensen_dice")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.90463992974239
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 22.74041276346604

$ws = $wb.Worksheets.Item("Dice")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.90463992974239
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 22.74041276346604

$ws = $wb.Worksheets.Item("Humman")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 26.88085480093671
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 75.98433840749414

$ws = $wb.Worksheets.Item("Wong2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 26.88085480093671
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 75.98433840749414

$ws = $wb.Worksheets.Item("Euclid")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 26.88085480093671
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 75.98433840749414

$ws = $wb.Worksheets.Item("Zoltar")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 2.083943208430914
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 23.64607728337236

$ws = $wb.Worksheets.Item("Rogot2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 0.9989754098360663
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 40.27005269320841

$ws = $wb.Worksheets.Item("Hamming")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 26.88085480093671
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 75.98433840749414

$ws = $wb.Worksheets.Item("Fleiss")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 11.84865339578455
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 66.63861241217795

$ws = $wb.Worksheets.Item("Anderberg")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.90463992974239
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 22.74041276346604

$ws = $wb.Worksheets.Item("Goodman")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.90463992974239
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 22.74041276346604

$ws = $wb.Worksheets.Item("Harmonic_mean")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 0.9989754098360663
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 41.53066451990631

$ws = $wb.Worksheets.Item("Kulczynski2")
$ws.Range("A14").Value = "Best exam"
$ws.Range("C14").Value = 1.51492974238876
$ws.Range("A15").Value = "Worst exam"
$ws.Range("C15").Value = 24.86460772833722
